# Slide 2 (sldId 256 / cId 4267678481, the "Project PP" landing-page slide)
# gets a new "Geographic filter?" rectangle, matching the existing
# accent4-styled filter rectangles already on the slide (e.g. "Slider -
# Aroma"). The cleanest way to reproduce PowerPoint's theme "Shape Style"
# (lnRef/fillRef/effectRef/fontRef) via COM automation -- which has no
# direct setter for it -- is to duplicate a shape that already has that
# style, then move it into place and retarget its text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "Rectangle 12" / "Slider - Aroma" already carries the accent4 shape style
# we need to clone.
$template = $s.Shapes.Item(4)

$dupRange = $template.Duplicate()
$shp = $dupRange.Item(1)

$shp.Name = "Rectangle 1"

# EMU-exact placement (465085, 11385330, 2017986, 562304 EMU).
$shp.Left = 36.62086684173229
$shp.Top = 896.4826965653543
$shp.Width = 158.89653783307085
$shp.Height = 44.275905511811025

$shp.TextFrame.TextRange.Text = "Geographic filter?"
